$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Find-ParagraphByText($doc, $target) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -eq $target) {
            return $p
        }
    }
    return $null
}

# --- 1. Replace the lone page-break paragraph with an empty spacer paragraph
#        plus a new paragraph announcing the VS2017 / June 2021 GDK update. ---
$cr = [string][char]13
$ff = [string][char]12
$pageBreakText = $ff + $cr
$pPageBreak = Find-ParagraphByText $d $pageBreakText
$xmlPageBreak = "<w:p $wns><w:pPr><w:spacing w:after=`"160`" w:line=`"259`" w:lineRule=`"auto`"/></w:pPr></w:p>" +
                "<w:p $wns><w:pPr><w:spacing w:after=`"160`" w:line=`"259`" w:lineRule=`"auto`"/></w:pPr>" +
                "<w:r><w:t>Update: this sample now works with Visual Studio 2017 (15.9.38) using June 2021 GDK.</w:t></w:r></w:p>"
$pPageBreak.Range.InsertXML($xmlPageBreak) | Out-Null

# --- 2. Drop the stray lastRenderedPageBreak from the "Using the sample" heading ---
$usingSampleText = "Using the sample" + $cr
$pUsingSample = Find-ParagraphByText $d $usingSampleText
$xmlUsingSample = "<w:p $wns><w:pPr><w:pStyle w:val=`"Heading1`"/></w:pPr><w:r><w:t>Using the sample</w:t></w:r></w:p>"
$pUsingSample.Range.InsertXML($xmlUsingSample) | Out-Null

# --- 3. Add the lastRenderedPageBreak to the "Sample Start Screen" heading instead ---
$sampleStartText = "Sample Start Screen" + $cr
$pSampleStart = Find-ParagraphByText $d $sampleStartText
$xmlSampleStart = "<w:p $wns><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Sample Start Screen</w:t></w:r></w:p>"
$pSampleStart.Range.InsertXML($xmlSampleStart) | Out-Null

# --- 4. Fill in the final (blank) row of the "Update history" table ---
$historyTable = $null
for ($t = 1; $t -le $d.Tables.Count; $t++) {
    $tbl = $d.Tables.Item($t)
    if ($tbl.Cell(1, 1).Range.Text -like "Description*") {
        $historyTable = $tbl
    }
}
$lastRow = $historyTable.Rows.Count
$historyTable.Cell($lastRow, 1).Range.Paragraphs.Item(1).Range.InsertAfter("Updated to make compatible with VS2017") | Out-Null
$historyTable.Cell($lastRow, 2).Range.Paragraphs.Item(1).Range.InsertAfter("September 2020") | Out-Null
$historyTable.Cell($lastRow, 3).Range.Paragraphs.Item(1).Range.InsertAfter("1.1") | Out-Null

Write-Output "Edits applied."
